# Apply cell value updates per the commit diff (Sun Feb 12 23:18:20 UTC 2023 crypto price refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'318.19"
$ws.Range("E2").Value = "'2.53%"

$ws.Range("D3").Value = "'41.06"
$ws.Range("E3").Value = "'-0.91%"

$ws.Range("D4").Value = "'5.135"
$ws.Range("E4").Value = "'0.09%"

$ws.Range("D5").Value = "'0.07641"
$ws.Range("E5").Value = "'-0.64%"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.692"
$ws.Range("E6").Value = "'4.35%"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9361"
$ws.Range("E7").Value = "'1.44%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.425"
$ws.Range("E8").Value = "'-1.35%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1243"
$ws.Range("E9").Value = "'1.80%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1826"
$ws.Range("E10").Value = "'-0.28%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09068"
$ws.Range("E11").Value = "'-1.38%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04182"
$ws.Range("E12").Value = "'-3.25%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1057"
$ws.Range("E13").Value = "'0.79%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001271"
$ws.Range("E14").Value = "'1.79%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005788"
$ws.Range("E15").Value = "'-0.46%"

$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007491"
$ws.Range("E16").Value = "'1,897.31%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.352"
$ws.Range("E17").Value = "'0.18%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.336"
$ws.Range("E18").Value = "'0.68%"

$ws.Range("D20").Value = "'8.418"
$ws.Range("E20").Value = "'21.22%"

$ws.Range("D21").Value = "'0.1349"
$ws.Range("E21").Value = "'-2.91%"

$ws.Range("D22").Value = "'0.2738"
$ws.Range("E22").Value = "'2.23%"

$ws.Range("D23").Value = "'0.04053"
$ws.Range("E23").Value = "'0.02%"

$ws.Range("D24").Value = "'0.001268"
$ws.Range("E24").Value = "'0.36%"

$ws.Range("D25").Value = "'0.004076"
$ws.Range("E25").Value = "'-0.56%"

$ws.Range("D26").Value = "'0.0001274"
$ws.Range("E26").Value = "'0.46%"

$ws.Range("D38").Value = "'0.02496"
$ws.Range("E38").Value = "'1.14%"

$ws.Range("D39").Value = "'0.05233"
$ws.Range("E39").Value = "'-0.68%"

$ws.Range("D40").Value = "'0.007769"
$ws.Range("E40").Value = "'-0.70%"

$ws.Range("D41").Value = "'0.1299"
$ws.Range("E41").Value = "'-1.20%"

$ws.Range("D42").Value = "'0.007242"
$ws.Range("E42").Value = "'6.64%"

$ws.Range("D43").Value = "'0.002171"
$ws.Range("E43").Value = "'17.95%"

$ws.Range("D44").Value = "'0.008245"
$ws.Range("E44").Value = "'0.87%"

$ws.Range("D45").Value = "'0.3180"
$ws.Range("E45").Value = "'2.50%"

$ws.Range("D46").Value = "'0.00006654"
$ws.Range("E46").Value = "'-0.90%"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.46%"

$ws.Range("D48").Value = "'0.2439"
$ws.Range("E48").Value = "'18.69%"

$ws.Range("D49").Value = "'0.004222"
$ws.Range("E49").Value = "'3.14%"

$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.46%"

$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.46%"
